$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped from
# 45204 (2023-10-05) to 45207 (2023-10-08) for every data row (rows 2-173).
$oldValue = 45204
$newValue = 45207

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
